$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data table. Insert a new row
# at row 370 (shifting the existing rows 370-457 down to 371-458) and
# populate it with the new record's values.
$ws.Rows.Item(370).Insert()

$ws.Cells.Item(370, 1).Value = 4
$ws.Cells.Item(370, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(370, 3).Value = "Los Lagos"
$ws.Cells.Item(370, 4).Value = 45173
$ws.Cells.Item(370, 5).Value = 10
$ws.Cells.Item(370, 6).Value = 100112021
$ws.Cells.Item(370, 7).Value = "Ají"
$ws.Cells.Item(370, 8).Value = "Inferno"
$ws.Cells.Item(370, 9).Value = "Primera"
$ws.Cells.Item(370, 10).Value = 70
$ws.Cells.Item(370, 11).Value = 29000
$ws.Cells.Item(370, 12).Value = 29000
$ws.Cells.Item(370, 13).Value = 29000
$ws.Cells.Item(370, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(370, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(370, 16).Value = 2900
$ws.Cells.Item(370, 17).Value = 10
$ws.Cells.Item(370, 18).Value = "Hortaliza"
